# Update "想去人数" (want-to-go count) values in column F
# for worksheets "展览" and "全部类型".

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 94
$ws1.Range("F4").Value = 70
$ws1.Range("F5").Value = 2495
$ws1.Range("F6").Value = 229
$ws1.Range("F7").Value = 376

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 94
$ws4.Range("F4").Value = 70
$ws4.Range("F5").Value = 2495
$ws4.Range("F6").Value = 229
$ws4.Range("F9").Value = 376
